# Adiciona novas pessoas a cena: duas novas linhas de dados na planilha de rotina
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(4, 4, 12, 0, 6, 28, 13, 0, 23, 0, 10, 0),
    @(3, 3, 13, 0, 12, 0, 12, 0, 22, 0, 9, 0)
)

$startRow = 4
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $rowValues.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
}

[void]$ws.Range("J6").Select()
